$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New combined values for rows 2-9 (card name + attribute list collapsed into
# a single Python-tuple-like text string)
$ws.Range("A2").Value = "('Hellspark Elemental', ['{1}{R}', 'Creature — Elemental', 'Trample, haste', 'At the beginning of the end step, sacrifice Hellspark Elemental.', 'Unearth {1}{R} ({1}{R}: Return this card from your graveyard to the battlefield. It gains haste. Exile it at the beginning of the next end step or if it would leave the battlefield. Unearth only as a sorcery.)', '3/1'])"
$ws.Range("A3").Value = "('Kor Duelist', ['{W}', 'Creature — Kor Soldier', 'As long as Kor Duelist is equipped, it has double strike. (It deals both first-strike and regular combat damage.)', '1/1'])"
$ws.Range("A4").Value = "(`"Marisi's Twinclaws`", ['{2}{R/W}{G}', 'Creature — Cat Warrior', 'Double strike', '2/4'])"
$ws.Range("A5").Value = "('Mind Control', ['{3}{U}{U}', 'Enchantment — Aura', 'Enchant creature', 'You control enchanted creature.'])"
$ws.Range("A6").Value = "('Path to Exile', ['{W}', 'Instant', 'Exile target creature. Its controller may search their library for a basic land card, put that card onto the battlefield tapped, then shuffle their library.'])"
$ws.Range("A7").Value = "('Rise from the Grave', ['{4}{B}', 'Sorcery', 'Put target creature card from a graveyard onto the battlefield under your control. That creature is a black Zombie in addition to its other colors and types.'])"
$ws.Range("A8").Value = "('Slave of Bolas', ['{3}{U/R}{B}', 'Sorcery', 'Gain control of target creature. Untap that creature. It gains haste until end of turn. Sacrifice it at the beginning of the next end step.'])"
$ws.Range("A9").Value = "('Vampire Nighthawk', ['{1}{B}{B}', 'Creature — Vampire Shaman', 'Flying, deathtouch, lifelink', '2/3'])"

# Remove the now-unused rows 10-40 entirely so the sheet's used range
# shrinks back down to A1:A9
$ws.Range("A10:A40").EntireRow.Delete()
